$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 8 (year 2025) values per the diff
$ws.Range("C8").Value = 1219
$ws.Range("E8").Value = 1021
$ws.Range("G8").Value = 83.7571780147662
$ws.Range("H8").Value = 16.2428219852338
